$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Update the L20 input value (0.106 -> 0.176); dependent formulas recalc automatically
$ws.Range("L20").Value = 0.176

# Update the active view: scroll position and selected cell
$ws.Activate()
$ws.Range("L6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
